$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 302.9375
$ws.Range("I2").Value = 310.5
$ws.Range("K2").Value = 310.5
$ws.Range("M2").Value = -197.5

$ws.Range("H138").Value = 2653.6191
$ws.Range("J138").Value = 3305.7273
$ws.Range("L138").Value = 9917.1819
$ws.Range("N138").Value = -20197.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1302.0385
$ws.Range("I2").Value = 1014.44446
$ws.Range("J2").Value = 1949.125
$ws.Range("K2").Value = 1014.44446
$ws.Range("L2").Value = 1949.125
$ws.Range("M2").Value = -901.44446
$ws.Range("N2").Value = -2175.125

$ws.Range("H4").Value = 758.8
$ws.Range("I4").Value = 758.8
$ws.Range("K4").Value = 758.8
$ws.Range("M4").Value = -642.8

$ws.Range("H5").Value = 437.54544
$ws.Range("I5").Value = 328.3
$ws.Range("K5").Value = 328.3
$ws.Range("M5").Value = -216.3

$ws.Range("H45").Value = 4793.7
$ws.Range("I45").Value = 4879.625
$ws.Range("J45").Value = 4450
$ws.Range("K45").Value = 4879.625
$ws.Range("L45").Value = 4450
$ws.Range("M45").Value = -4502.625
$ws.Range("N45").Value = -5204

$ws.Range("H97").Value = 1985.25
$ws.Range("I97").Value = 1256.6364
$ws.Range("K97").Value = 1256.6364
$ws.Range("M97").Value = -760.6364000000001

$ws.Range("H116").Value = 1302.0385
$ws.Range("I116").Value = 1014.44446
$ws.Range("J116").Value = 1949.125
$ws.Range("K116").Value = 1014.44446
$ws.Range("L116").Value = 1949.125
$ws.Range("M116").Value = 1279.55554
$ws.Range("N116").Value = -6537.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1302.0385
$ws.Range("I3").Value = 1014.44446
$ws.Range("J3").Value = 1949.125
$ws.Range("K3").Value = 1014.44446
$ws.Range("L3").Value = 1949.125
$ws.Range("M3").Value = -900.44446
$ws.Range("N3").Value = -2177.125

$ws.Range("H4").Value = 437.54544
$ws.Range("I4").Value = 328.3
$ws.Range("K4").Value = 328.3
$ws.Range("M4").Value = -213.3

$ws.Range("H22").Value = 289.75
$ws.Range("I22").Value = 289.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 289.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -116.75
$ws.Range("N22").ClearContents()

$ws.Range("H38").Value = 23999.666
$ws.Range("J38").Value = 23999.666
$ws.Range("L38").Value = 23999.666
$ws.Range("N38").Value = -24831.666

$ws.Range("H105").Value = 3249.2856
$ws.Range("I105").Value = 3224.1667
$ws.Range("J105").Value = 3400
$ws.Range("K105").Value = 3224.1667
$ws.Range("L105").Value = 3400
$ws.Range("M105").Value = -1477.1667
$ws.Range("N105").Value = -6894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1156
$ws.Range("I22").Value = 207
$ws.Range("J22").Value = 1725.4
$ws.Range("K22").Value = 207
$ws.Range("L22").Value = 1725.4
$ws.Range("M22").Value = 143
$ws.Range("N22").Value = -2425.4

$ws.Range("H132").Value = 7783.3
$ws.Range("I132").Value = 7536.8887
$ws.Range("K132").Value = 22610.6661
$ws.Range("M132").Value = -20080.6661

$ws.Range("H134").Value = 1694.8
$ws.Range("I134").Value = 1694.8
$ws.Range("K134").Value = 5084.4
$ws.Range("M134").Value = -2549.4

$ws.Range("H141").Value = 401254.66
$ws.Range("J141").Value = 401254.66
$ws.Range("L141").Value = 401254.66
$ws.Range("N141").Value = -411614.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3400
$ws.Range("I3").Value = 3400
$ws.Range("K3").Value = 10200
$ws.Range("M3").Value = -10088

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61118

$ws.Range("H61").Value = 285.5
$ws.Range("I61").Value = 183.42857
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 550.28571
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -335.28571
$ws.Range("N61").Value = -3430

$ws.Range("H62").Value = 9366.166999999999
$ws.Range("J62").Value = 9366.166999999999
$ws.Range("L62").Value = 28098.501
$ws.Range("N62").Value = -29470.501

$ws.Range("H65").Value = 9366.166999999999
$ws.Range("J65").Value = 9366.166999999999
$ws.Range("L65").Value = 84295.503
$ws.Range("N65").Value = -91159.503

$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -16872

$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -54360

$ws.Range("H94").Value = 10214
$ws.Range("I94").Value = 749
$ws.Range("K94").Value = 2247
$ws.Range("M94").Value = -1571

$ws.Range("H106").Value = 17749.75
$ws.Range("J106").Value = 19999.834
$ws.Range("L106").Value = 59999.50199999999
$ws.Range("N106").Value = -61891.50199999999

$ws.Range("H107").Value = 133
$ws.Range("I107").Value = 111.25
$ws.Range("J107").Value = 167.8
$ws.Range("K107").Value = 333.75
$ws.Range("L107").Value = 503.4
$ws.Range("M107").Value = 1586.25
$ws.Range("N107").Value = -4343.4

$ws.Range("H114").Value = 800
$ws.Range("I114").Value = 800
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 2400
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 854
$ws.Range("N114").ClearContents()

$ws.Range("H137").Value = 35333.668
$ws.Range("J137").Value = 3000.5
$ws.Range("L137").Value = 9001.5
$ws.Range("N137").Value = -19201.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7999.5
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 7999
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 7999
$ws.Range("M41").Value = -7645
$ws.Range("N41").Value = -8709

$ws.Range("H80").Value = 10566.583
$ws.Range("I80").Value = 9000
$ws.Range("K80").Value = 9000
$ws.Range("M80").Value = -8002

$ws.Range("H83").Value = 10566.583
$ws.Range("I83").Value = 9000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40008

$ws.Range("H102").Value = 1845.1428
$ws.Range("I102").Value = 1394.1305
$ws.Range("K102").Value = 1394.1305
$ws.Range("M102").Value = 227.8695

$ws.Range("H107").Value = 2254.3333
$ws.Range("I107").Value = 716
$ws.Range("J107").Value = 5331
$ws.Range("K107").Value = 716
$ws.Range("L107").Value = 5331
$ws.Range("M107").Value = 1204
$ws.Range("N107").Value = -9171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 940.6
$ws.Range("I16").Value = 940.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 940.6
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -770.6
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 1489.25
$ws.Range("I22").Value = 1477
$ws.Range("J22").Value = 1491
$ws.Range("K22").Value = 1477
$ws.Range("L22").Value = 1491
$ws.Range("M22").Value = -1182
$ws.Range("N22").Value = -2081

$ws.Range("H27").Value = 1489.25
$ws.Range("I27").Value = 1477
$ws.Range("J27").Value = 1491
$ws.Range("K27").Value = 1477
$ws.Range("L27").Value = 1491
$ws.Range("M27").Value = -1370
$ws.Range("N27").Value = -1705

$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3864

$ws.Range("H46").Value = 4249.25
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 4999
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 4999
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -5375

$ws.Range("H62").Value = 16000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 34999
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H65").Value = 16000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 34999
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 44553.875
$ws.Range("I136").Value = 44915
$ws.Range("K136").Value = 134745
$ws.Range("M136").Value = -132195

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
